# Auto-generated Excel COM-interop edit script
# Updates crypto price/volume data per the commit diff (Tue Jul 18 08:11:46 UTC 2023 GitHub Actions run)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Price" column (D): values are textual (dotted/grouped numbers), so force
# each cell to Text format before writing to avoid Excel auto-converting the
# string into a numeric value.
$priceUpdates = [ordered]@{
    'D2' = '29.985.07'
    'D3' = '1.898.79'
    'D5' = '0.7406'
    'D6' = '241.47'
    'D8' = '0.3059'
    'D9' = '25.73'
    'D10' = '0.06891'
    'D11' = '0.08015'
    'D12' = '0.7559'
    'D13' = '1.903.36'
    'D14' = '5.220'
    'D15' = '91.03'
    'D16' = '6.165'
    'D17' = '29.995.18'
    'D18' = '13.98'
    'D19' = '0.000007746'
    'D22' = '2.159.26'
    'D24' = '7.065'
    'D25' = '167.31'
    'D26' = '9.292'
    'D27' = '18.80'
    'D28' = '0.1261'
    'D29' = '2.027'
    'D31' = '1.527'
    'D32' = '4.290'
    'D33' = '4.030'
    'D34' = '0.05255'
    'D35' = '1.282'
    'D36' = '0.7370'
    'D38' = '0.01935'
    'D39' = '2.771'
    'D40' = '6.249'
    'D41' = '0.4439'
    'D42' = '72.60'
    'D43' = '1.951'
    'D45' = '0.8321'
    'D46' = '7.657'
    'D47' = '101.01'
    'D48' = '9.813'
    'D49' = '2.051.66'
    'D50' = '36.47'
    'D51' = '0.1162'
}
foreach ($ref in $priceUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $priceUpdates[$ref]
}

# --- "Volume(1h)" column (E), plus the Coin-name / Link swap for rows 25 & 26
# (Cosmos <-> Monero exchanged ranking positions).
$otherUpdates = [ordered]@{
    'E2' = '  -1.19%  '
    'E3' = '  -1.75%  '
    'E4' = '  +0.04%  '
    'E5' = '  -0.90%  '
    'E6' = '  -0.67%  '
    'E8' = '  -3.33%  '
    'E9' = '  -6.76%  '
    'E10' = '  -3.37%  '
    'E11' = '  -0.64%  '
    'E12' = '  -2.92%  '
    'E13' = '  -1.34%  '
    'E14' = '  -3.30%  '
    'E15' = '  -2.26%  '
    'E16' = '  +2.22%  '
    'E17' = '  -1.07%  '
    'E18' = '  -4.00%  '
    'E19' = '  -2.20%  '
    'E20' = '  -5.85%  '
    'E22' = '  -0.21%  '
    'E23' = '  +0.08%  '
    'E24' = '  +5.92%  '
    'B25' = 'Monero'
    'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E25' = '  +0.93%  '
    'B26' = 'Cosmos'
    'C26' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E26' = '  -2.78%  '
    'E27' = '  -1.42%  '
    'E28' = '  -2.65%  '
    'E29' = '  -7.07%  '
    'E30' = '  -2.04%  '
    'E31' = '  -1.84%  '
    'E32' = '  -2.85%  '
    'E33' = '  -2.78%  '
    'E34' = '  +0.26%  '
    'E35' = '  -2.53%  '
    'E36' = '  -2.51%  '
    'E37' = '  -2.06%  '
    'E38' = '  -0.56%  '
    'E39' = '  -1.10%  '
    'E40' = '  -4.00%  '
    'E41' = '  -1.94%  '
    'E42' = '  -6.96%  '
    'E43' = '  -1.15%  '
    'E44' = '  -0.07%  '
    'E45' = '  -1.26%  '
    'E46' = '  -0.22%  '
    'E47' = '  -0.73%  '
    'E48' = '  -1.47%  '
    'E49' = '  -0.89%  '
    'E50' = '  -3.56%  '
    'E51' = '  -4.98%  '
}
foreach ($ref in $otherUpdates.Keys) {
    $ws.Range($ref).Value = $otherUpdates[$ref]
}

